$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O1").Value = "objetivos_relacionados"
$ws.Range("O1").Select()
